$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 138914.46882130537
$ws.Range("D2").Value = 20.5

$ws.Range("C3").Value = 14516.597729148474
$ws.Range("D3").Value = 11.4

$ws.Range("C4").Value = 9845.141984185488
$ws.Range("D4").Value = 15.3

$ws.Range("C5").Value = 114552.72910797142
$ws.Range("D5").Value = 23.5

$ws.Range("C6").Value = 26125.173186287033
$ws.Range("D6").Value = 21.5

$ws.Range("C7").Value = 33849.03757102444
$ws.Range("D7").Value = 19.8

$ws.Range("C8").Value = 39263.165425997526
$ws.Range("D8").Value = 19.5

$ws.Range("C9").Value = 6052.406386566344
$ws.Range("D9").Value = 24.0

$ws.Range("C10").Value = 33595.9699673872
$ws.Range("D10").Value = 20.9

$ws.Range("C11").Value = 28.716284042823688
$ws.Range("D11").Value = 22.0
